$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells F393:G441 per diff
$ws.Range("F393").Value = 308264
$ws.Range("G393").Value = 1231
$ws.Range("F394").Value = 166133
$ws.Range("F395").Value = 752816
$ws.Range("G395").Value = 1957
$ws.Range("F396").Value = 164839
$ws.Range("F397").Value = 107838
$ws.Range("F398").Value = 298877
$ws.Range("G398").Value = 1467
$ws.Range("F399").Value = 201127
$ws.Range("G399").Value = 964
$ws.Range("F400").Value = 150191
$ws.Range("F401").Value = 272129
$ws.Range("G401").Value = 929
$ws.Range("F402").Value = 719603
$ws.Range("F403").Value = 351843
$ws.Range("F404").Value = 225015
$ws.Range("G404").Value = 913
$ws.Range("F405").Value = 173584
$ws.Range("G405").Value = 691
$ws.Range("F406").Value = 170374
$ws.Range("G406").Value = 680
$ws.Range("F407").Value = 157889
$ws.Range("G407").Value = 679
$ws.Range("F408").Value = 303834
$ws.Range("F409").Value = 706207
$ws.Range("F411").Value = 224899
$ws.Range("G411").Value = 826
$ws.Range("F412").Value = 175874
$ws.Range("G412").Value = 644
$ws.Range("F413").Value = 149121
$ws.Range("G413").Value = 658
$ws.Range("F414").Value = 148587
$ws.Range("G414").Value = 563
$ws.Range("F415").Value = 306548
$ws.Range("F416").Value = 669421
$ws.Range("F417").Value = 340904
$ws.Range("F418").Value = 202351
$ws.Range("G418").Value = 703
$ws.Range("F419").Value = 148985
$ws.Range("G419").Value = 510
$ws.Range("F420").Value = 138245
$ws.Range("F421").Value = 152691
$ws.Range("G421").Value = 532
$ws.Range("F422").Value = 297203
$ws.Range("G422").Value = 644
$ws.Range("F423").Value = 436055
$ws.Range("F424").Value = 264608
$ws.Range("F425").Value = 139138
$ws.Range("F429").Value = 177892
$ws.Range("F430").Value = 173579
$ws.Range("F431").Value = 170485
$ws.Range("F432").Value = 123128
$ws.Range("F434").Value = 79640
$ws.Range("F435").Value = 81865
$ws.Range("F436").Value = 143266
$ws.Range("F437").Value = 165018
$ws.Range("F438").Value = 120690
$ws.Range("F439").Value = 88009
$ws.Range("F440").Value = 72955
$ws.Range("G440").Value = 223
$ws.Range("F441").Value = 65967
$ws.Range("G441").Value = 200

# Add new row 442
$ws.Range("A442").Value = 44336
$ws.Range("A442").NumberFormat = "yyyy-mm-dd"
$ws.Range("B442").Value = 388529
$ws.Range("C442").Value = 4004
$ws.Range("D442").Value = 138
$ws.Range("E442").Value = 12280
$ws.Range("F442").Value = 54523
$ws.Range("G442").Value = 137
